# Bulk sample typist and typist_qc added and also in manual
#
# Effect on sheet "with_all_correctdata" (A1:M3):
#   - Columns E,F,G are left-rotated: new E = old F, new F = old G, new G = old E
#       (Client/Typist/Typist QC  ->  Typist/Typist QC/Client)
#   - Columns H and J are swapped: new H = old J, new J = old H
#       (Product Name/Lob/Process  ->  Lob/Process/Product Name order change)
#   - Column I is left untouched.
# This is applied uniformly to the header row (1) and the two data rows (2-3).
#
# We stage the source values (together with their existing cell formatting) in
# a scratch area of the sheet (columns far to the right, P:T) using
# Range.Copy(destination) -- this preserves both the literal value and the
# style index of each cell -- then copy the staged cells back into their
# final destinations and clear the scratch area again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 1, 2, 3

foreach ($r in $rows) {
    # Stage old E,F,G,H,J (value + formatting) into scratch columns P,Q,R,S,T
    $ws.Range("E$r").Copy($ws.Range("P$r")) | Out-Null
    $ws.Range("F$r").Copy($ws.Range("Q$r")) | Out-Null
    $ws.Range("G$r").Copy($ws.Range("R$r")) | Out-Null
    $ws.Range("H$r").Copy($ws.Range("S$r")) | Out-Null
    $ws.Range("J$r").Copy($ws.Range("T$r")) | Out-Null
}

foreach ($r in $rows) {
    # Rotate E,F,G left: E<-oldF(Q), F<-oldG(R), G<-oldE(P)
    $ws.Range("Q$r").Copy($ws.Range("E$r")) | Out-Null
    $ws.Range("R$r").Copy($ws.Range("F$r")) | Out-Null
    $ws.Range("P$r").Copy($ws.Range("G$r")) | Out-Null

    # Swap H and J: H<-oldJ(T), J<-oldH(S)
    $ws.Range("T$r").Copy($ws.Range("H$r")) | Out-Null
    $ws.Range("S$r").Copy($ws.Range("J$r")) | Out-Null
}

# Clear the scratch area used for staging
$ws.Range("P1:T3").Clear() | Out-Null

# The explicit column-width override that used to sit on column E now sits on
# column G (the <col min="5".../> entry moves to <col min="7".../>).
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# Update the active cell selection as recorded in the saved file.
$ws.Range("C15").Select() | Out-Null
